# update models and data
#
# Adds a new one-joint hip muscle ("piri" / piriformis) to the Muscles sheet
# and a new "Combined" column (column F) on every sheet that lists every
# unique value used in that sheet's row data.

$wb = $excel.ActiveWorkbook

$wsMuscles = $wb.Worksheets.Item("Muscles")
$wsDof     = $wb.Worksheets.Item("DegreesOfFreedom")
$wsScale   = $wb.Worksheets.Item("ScaleFactors")
$wsForces  = $wb.Worksheets.Item("ForcesOnBodies")

# --- Muscles: new "Combined" column F = full muscle list ---
$musclesCombined = @(
    "Combined",
    "bflh", "grac", "recfem", "sart", "semimem", "semiten", "tfl",
    "addlong", "addbrev", "addmagDist", "addmagMid", "addmagProx", "addmagIsch",
    "glmax1", "glmax2", "glmax3",
    "glmed1", "glmed2", "glmed3",
    "glmin1", "glmin2", "glmin3",
    "iliacus", "psoas", "piri",
    "vasmed", "vaslat", "vasint",
    "bfsh", "gaslat", "gasmed",
    "perbrev", "perlong",
    "tibant", "tibpost", "soleus"
)
for ($i = 0; $i -lt $musclesCombined.Length; $i++) {
    $row = $i + 1
    $wsMuscles.Cells.Item($row, 6).Value = $musclesCombined[$i]
}

# --- Muscles: add the new one-joint muscle "piri" under column B ---
$wsMuscles.Range("B19").Value = "piri"

# --- DegreesOfFreedom: new "Combined" column F ---
$dofCombined = @("Combined", "HF", "HA", "HR", "KF", "AF", "ST")
for ($i = 0; $i -lt $dofCombined.Length; $i++) {
    $row = $i + 1
    $wsDof.Cells.Item($row, 6).Value = $dofCombined[$i]
}

# --- ScaleFactors: new "Combined" column F ---
$scaleCombined = @("Combined", "P1", "P2", "P3", "Fem", "Tib", "Foot")
for ($i = 0; $i -lt $scaleCombined.Length; $i++) {
    $row = $i + 1
    $wsScale.Cells.Item($row, 6).Value = $scaleCombined[$i]
}

# --- ForcesOnBodies: new "Combined" column F ---
$forcesCombined = @("Combined", "pelvis", "femur", "tibia", "calc")
for ($i = 0; $i -lt $forcesCombined.Length; $i++) {
    $row = $i + 1
    $wsForces.Cells.Item($row, 6).Value = $forcesCombined[$i]
}

# --- Restore per-sheet selections left behind by the editor ---
$wsMuscles.Range("P7").Select()
$wsDof.Range("F10").Select()
$wsScale.Range("G12").Select()
$wsForces.Range("F6").Select()

# ForcesOnBodies ends up the active/visible tab after the edit.
$wsForces.Activate()
